$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1782.7391
$ws.Range("I15").Value = 1782.7391
$ws.Range("K15").Value = 5348.2173
$ws.Range("M15").Value = -5179.2173
# Row 17
$ws.Range("H17").Value = 763617.9399999999
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 763617.9399999999
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2290853.82
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2291189.82
# Row 31
$ws.Range("H31").Value = 1237.625
$ws.Range("J31").Value = 3002
$ws.Range("L31").Value = 9006
$ws.Range("N31").Value = -9466
# Row 38
$ws.Range("H38").Value = 242.9
$ws.Range("I38").Value = 242.9
$ws.Range("K38").Value = 728.7
$ws.Range("M38").Value = -356.7
# Row 41
$ws.Range("H41").Value = 850.8570999999999
$ws.Range("I41").Value = 1267
$ws.Range("K41").Value = 1267
$ws.Range("M41").Value = -827
# Row 61
$ws.Range("H61").Value = 884.5
$ws.Range("I61").Value = 884.5
$ws.Range("K61").Value = 2653.5
$ws.Range("M61").Value = -2481.5
# Row 88
$ws.Range("H88").Value = 12368.667
$ws.Range("I88").Value = 9838
$ws.Range("J88").Value = 12685
$ws.Range("K88").Value = 9838
$ws.Range("L88").Value = 12685
$ws.Range("M88").Value = -9432
$ws.Range("N88").Value = -13497
# Row 91
$ws.Range("H91").Value = 12368.667
$ws.Range("I91").Value = 9838
$ws.Range("J91").Value = 12685
$ws.Range("K91").Value = 9838
$ws.Range("L91").Value = 12685
$ws.Range("M91").Value = -8434
$ws.Range("N91").Value = -15493
# Row 98
$ws.Range("H98").Value = 1416.52
$ws.Range("J98").Value = 5500
$ws.Range("L98").Value = 5500
$ws.Range("N98").Value = -8496
# Row 111
$ws.Range("H111").Value = 201905.4
$ws.Range("I111").Value = 1676
$ws.Range("J111").Value = 502249.5
$ws.Range("K111").Value = 5028
$ws.Range("L111").Value = 1506748.5
$ws.Range("M111").Value = -1961
$ws.Range("N111").Value = -1512882.5
# Row 112
$ws.Range("H112").Value = 103092.9
$ws.Range("I112").Value = 204829.8
$ws.Range("J112").Value = 69180.60000000001
$ws.Range("K112").Value = 614489.3999999999
$ws.Range("L112").Value = 207541.8
$ws.Range("M112").Value = -613381.3999999999
$ws.Range("N112").Value = -209757.8
# Row 122
$ws.Range("H122").Value = 1416.52
$ws.Range("J122").Value = 5500
$ws.Range("L122").Value = 16500
$ws.Range("N122").Value = -21400

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3562.125
$ws.Range("I32").Value = 2170.6843
$ws.Range("K32").Value = 2170.6843
$ws.Range("M32").Value = -1883.6843
# Row 45
$ws.Range("H45").Value = 5291.913
$ws.Range("I45").Value = 7100.6665
$ws.Range("K45").Value = 7100.6665
$ws.Range("M45").Value = -6723.6665
# Row 63
$ws.Range("H63").Value = 999.5
$ws.Range("I63").Value = 999.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 999.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -313.5
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 999.5
$ws.Range("I66").Value = 999.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 4997.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -1565.5
$ws.Range("N66").ClearContents()
# Row 110
$ws.Range("H110").Value = 59980.707
$ws.Range("I110").Value = 72648.21000000001
$ws.Range("K110").Value = 72648.21000000001
$ws.Range("M110").Value = -70603.21000000001
# Row 123
$ws.Range("H123").Value = 49995
$ws.Range("J123").Value = 49995
$ws.Range("L123").Value = 49995
$ws.Range("N123").Value = -59795

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 3219.7
$ws.Range("I22").Value = 4640.4
$ws.Range("K22").Value = 4640.4
$ws.Range("M22").Value = -4467.4

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 1669
$ws.Range("I5").Value = 1651.75
$ws.Range("J5").Value = 1703.5
$ws.Range("K5").Value = 1651.75
$ws.Range("L5").Value = 1703.5
$ws.Range("M5").Value = -1539.75
$ws.Range("N5").Value = -1927.5
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 86
$ws.Range("H86").Value = 11585.182
$ws.Range("I86").Value = 9078.299999999999
$ws.Range("K86").Value = 9078.299999999999
$ws.Range("M86").Value = -7955.299999999999
# Row 89
$ws.Range("H89").Value = 11585.182
$ws.Range("I89").Value = 9078.299999999999
$ws.Range("K89").Value = 45391.5
$ws.Range("M89").Value = -39775.5
# Row 132
$ws.Range("H132").Value = 71431736
$ws.Range("I132").Value = 90912210
$ws.Range("K132").Value = 272736630
$ws.Range("M132").Value = -272734100

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 167886.38
$ws.Range("J11").Value = 50001
$ws.Range("L11").Value = 150003
$ws.Range("N11").Value = -150283
# Row 17
$ws.Range("H17").Value = 2572.818
$ws.Range("I17").Value = 2928.4285
$ws.Range("K17").Value = 8785.2855
$ws.Range("M17").Value = -8616.2855
# Row 50
$ws.Range("H50").Value = 1984.6364
$ws.Range("I50").Value = 1544.7142
$ws.Range("K50").Value = 4634.142599999999
$ws.Range("M50").Value = -4153.142599999999
# Row 53
$ws.Range("H53").Value = 1984.6364
$ws.Range("I53").Value = 1544.7142
$ws.Range("K53").Value = 4634.142599999999
$ws.Range("M53").Value = -4153.142599999999
# Row 69
$ws.Range("H69").Value = 825
$ws.Range("J69").Value = 900
$ws.Range("L69").Value = 2700
$ws.Range("N69").Value = -4322
# Row 72
$ws.Range("H72").Value = 825
$ws.Range("J72").Value = 900
$ws.Range("L72").Value = 8100
$ws.Range("N72").Value = -16212
# Row 131
$ws.Range("H131").Value = 1709.9375
$ws.Range("I131").Value = 1054.7858
$ws.Range("K131").Value = 3164.3574
$ws.Range("M131").Value = 1875.6426
# Row 137
$ws.Range("H137").Value = 9093670
$ws.Range("J137").Value = 3141.25
$ws.Range("L137").Value = 9423.75
$ws.Range("N137").Value = -19623.75
# Row 140
$ws.Range("H140").Value = 1662.4706
$ws.Range("I140").Value = 1662.4706
$ws.Range("K140").Value = 4987.4118
$ws.Range("M140").Value = 192.5882000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3193.75
$ws.Range("I80").Value = 3258.6667
$ws.Range("K80").Value = 3258.6667
$ws.Range("M80").Value = -2260.6667
# Row 83
$ws.Range("H83").Value = 3193.75
$ws.Range("I83").Value = 3258.6667
$ws.Range("K83").Value = 16293.3335
$ws.Range("M83").Value = -11301.3335
# Row 97
$ws.Range("H97").Value = 2601.28
$ws.Range("I97").Value = 2654.4119
$ws.Range("J97").Value = 2488.375
$ws.Range("K97").Value = 2654.4119
$ws.Range("L97").Value = 2488.375
$ws.Range("M97").Value = -2158.4119
$ws.Range("N97").Value = -3480.375
# Row 136
$ws.Range("H136").Value = 38750
$ws.Range("J136").Value = 38750
$ws.Range("L136").Value = 116250
$ws.Range("N136").Value = -121350

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1277.5714
$ws.Range("I46").Value = 1299
$ws.Range("K46").Value = 1299
$ws.Range("M46").Value = -1111
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
# Row 68
$ws.Range("H68").Value = 5265355.5
$ws.Range("I68").Value = 6580995.5
$ws.Range("K68").Value = 6580995.5
$ws.Range("M68").Value = -6580246.5
# Row 71
$ws.Range("H71").Value = 5265355.5
$ws.Range("I71").Value = 6580995.5
$ws.Range("K71").Value = 32904977.5
$ws.Range("M71").Value = -32901233.5
# Row 82
$ws.Range("H82").Value = 1182.0869
$ws.Range("I82").Value = 1213.3158
$ws.Range("K82").Value = 1213.3158
$ws.Range("M82").Value = -852.3158000000001
# Row 85
$ws.Range("H85").Value = 1182.0869
$ws.Range("I85").Value = 1213.3158
$ws.Range("K85").Value = 1213.3158
$ws.Range("M85").Value = 34.68419999999992

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 34243.25
$ws.Range("J41").Value = 34243.25
$ws.Range("L41").Value = 34243.25
$ws.Range("N41").Value = -35023.25
# Row 96
$ws.Range("H96").Value = 1582.75
$ws.Range("I96").Value = 1169.6666
$ws.Range("J96").Value = 1995.8334
$ws.Range("K96").Value = 1169.6666
$ws.Range("L96").Value = 1995.8334
$ws.Range("M96").Value = 203.3334
$ws.Range("N96").Value = -4741.8334
# Row 113
$ws.Range("H113").Value = 693.13336
$ws.Range("I113").Value = 638.96155
$ws.Range("K113").Value = 1916.88465
$ws.Range("M113").Value = 253.11535
# Row 122
$ws.Range("H122").Value = 2183.4285
$ws.Range("I122").Value = 2118.8
$ws.Range("K122").Value = 6356.400000000001
$ws.Range("M122").Value = -3906.400000000001
